# Update the acquisition timestamp (column A) for rows 2-23 on the
# "ランサーズ" sheet from "2025-11-20 18:26:42" to "2025-11-20 18:34:29".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$oldValue = "2025-11-20 18:26:42"
$newValue = "2025-11-20 18:34:29"

for ($row = 2; $row -le 23; $row++) {
    $cell = $ws.Cells.Item($row, 1)
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
    }
}
